$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive updated values, in order, for each data row (B:F then I:N, skipping G and H)
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

# New vm_pu values for rows 2-25 (one data row per source row), columns per $cols above
$data = @"
1.02,1.025460861056469,1.028105407604875,1.025790516327879,1.023960685446972,1.030388344066191,1.030629755629894,1.030923192599291,1.02861505614883,1.026790588204794,1.014159422784409
1.02,1.026867341249778,1.029424797213447,1.026996600517514,1.026010072726327,1.030683074228822,1.031673222656614,1.032049067396704,1.029627436518428,1.028643585896856,1.014511244074258
1.02,1.0277750418129,1.030276484227645,1.027775288599072,1.027333302738543,1.030870812285324,1.032345655912636,1.032775009804631,1.03028023224882,1.02983938402983,1.014737793203968
1.02,1.02815607688255,1.030634051374047,1.028102241899186,1.027888922542228,1.030949028163207,1.032627693080332,1.03307958698807,1.030554127557015,1.030341344766503,1.014832772720897
1.02,1.028220021636942,1.030694060427603,1.028157115115534,1.027982175033278,1.030962119431945,1.032675010223273,1.033130691332289,1.030600084299744,1.030425582527825,1.014848704916755
1.02,1.027780135420362,1.030281263942022,1.02777965895669,1.027340729551889,1.030871860193597,1.032349427067536,1.032779081964063,1.03028389416552,1.029846094191003,1.014739063350599
1.02,1.025936687143139,1.028551729923946,1.026198480229262,1.024653889107603,1.03048856629777,1.030982975784741,1.031304224543661,1.028957671172283,1.027417493441487,1.014278552145158
1.02,1.022669613712598,1.025488044742691,1.023398692957669,1.019896579310227,1.02979028524567,1.02855366593329,1.028685299552754,1.026602933319814,1.023112599173195,1.013458524035746
1.02,1.020478430348019,1.023434314755208,1.021522626627981,1.016708554406316,1.029309238037926,1.026919264336067,1.026925418105565,1.025020761336469,1.020224525301471,1.012905945977901
1.02,1.019526379405852,1.022542240738586,1.020707908557985,1.015323928820687,1.0290972219662,1.026207925872142,1.026159962093659,1.024332643874812,1.018969413376775,1.012665243155395
1.02,1.019172245155697,1.02221045493682,1.020404921554816,1.014808964775163,1.029017907853635,1.025943148749598,1.025875115416424,1.024076583740584,1.018502503472026,1.012575617263222
1.02,1.019248230968961,1.022281643725536,1.020469929925075,1.014919456207856,1.029034946477619,1.025999969567492,1.025936239774401,1.024131530574379,1.018602689491866,1.012594852239448
1.02,1.019497116826543,1.022514824045593,1.020682871041897,1.015281375175857,1.029090677317684,1.026186050671855,1.026136427284664,1.024311487349026,1.018930832970828,1.012657839118068
1.02,1.019650397005328,1.022658436889168,1.020814022543584,1.015504278248097,1.029124940376358,1.02630062766313,1.026259699942904,1.024422303108812,1.019132918891504,1.012696618419684
1.02,1.020541545316506,1.023493459039086,1.021576646146386,1.016800357216923,1.029323230218519,1.026966396268769,1.026976146164137,1.025066364982626,1.020307725165656,1.01292189019804
1.02,1.021099659996991,1.024016491077711,1.022054379558871,1.017612215451156,1.029446614223638,1.027383036879838,1.027424633104347,1.025469552094996,1.021043415469129,1.013062811594417
1.02,1.021424885642954,1.024321297773496,1.022332805386749,1.018085355938566,1.02951822336519,1.027625706479233,1.027685898998602,1.025704432841733,1.021472092923163,1.013144870605065
1.02,1.021535726396809,1.024425183640337,1.022427702984039,1.018246617006771,1.029542579489485,1.027708391467503,1.027774928393194,1.025784471906588,1.021618187148174,1.01317282727674
1.02,1.021039812051215,1.023960402589239,1.022003146928514,1.017525152563383,1.029433413407948,1.027338371531751,1.027376548790381,1.025426324139681,1.020964528342913,1.013047706362755
1.02,1.019423840042146,1.022446170197044,1.020620175319658,1.015174817257242,1.029074281510762,1.026131269797707,1.026077491541926,1.024258507400647,1.018834222537514,1.012639297097392
1.02,1.01840491497372,1.021491620592049,1.019748533358484,1.013693281688119,1.028845228459934,1.025369105363321,1.025257696642524,1.023521573776931,1.017490722695624,1.01238125040345
1.02,1.018945345045094,1.021997884883135,1.020210810399503,1.014479038385813,1.02896696316997,1.025773450516331,1.025692575318091,1.023912493037102,1.018203332279304,1.012518166579709
1.02,1.021066855729624,1.023985747394305,1.02202629745348,1.017564493778613,1.029439379396641,1.027358554960933,1.027398277038602,1.025445857892675,1.021000175393387,1.013054532193361
1.02,1.023516503063064,1.026282030831212,1.024124154796602,1.021129275707476,1.0299735330289,1.029184287537179,1.029364775925342,1.027213835790746,1.024228636120656,1.013671548825379
"@

$rowIndex = 2
foreach ($line in ($data -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $values = $line -split ","
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($rowIndex, $cols[$i]).Value = [double]$values[$i]
    }
    $rowIndex++
}

Write-Output "applied vm_pu updates"
